# Database schemas.xlsx — add "Project user" join table, a user_id FK on
# Projects, and a "comment" field on Task assignment.
# (commit: "added user roles, permissions, check permission of users")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Projects table: new FK column ---------------------------------------
$ws.Range("A17").Value = "user_id"

# --- Task assignment table: new "comment" field ---------------------------
$ws.Range("G18").Value = "comment"

# --- New "Project user" join table (rows 24-27, column A) -----------------
$ws.Range("A24").Value = "**Project user**"
$ws.Range("A25").Value = "id"
$ws.Range("A26").Value = "user_id"
$ws.Range("A27").Value = "project_id"

# --- Column A is now wider (longest label is "**Project user**"), so it is
# split off from the former A:D default-width group and auto-sized.
$ws.Columns.Item(1).ColumnWidth = 15.833333333333332

# --- Selection moved to E5 ---------------------------------------------
[void]$ws.Range("E5").Select()
